# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to the latest scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3046
$ws1.Range("F4").Value = 48
$ws1.Range("F5").Value = 35
$ws1.Range("F7").Value = 714
$ws1.Range("F8").Value = 14617
$ws1.Range("F9").Value = 170
$ws1.Range("F10").Value = 129
$ws1.Range("F11").Value = 5838
$ws1.Range("F12").Value = 595
$ws1.Range("F13").Value = 79
$ws1.Range("F14").Value = 45
$ws1.Range("F15").Value = 69
$ws1.Range("F17").Value = 16
$ws1.Range("F18").Value = 86
$ws1.Range("F19").Value = 186
$ws1.Range("F20").Value = 800
$ws1.Range("F22").Value = 61
$ws1.Range("F23").Value = 10623
$ws1.Range("F25").Value = 63
$ws1.Range("F26").Value = 90
$ws1.Range("F27").Value = 3739
$ws1.Range("F28").Value = 246
$ws1.Range("F29").Value = 68

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3046
$ws4.Range("F5").Value = 48
$ws4.Range("F6").Value = 35
$ws4.Range("F8").Value = 715
$ws4.Range("F9").Value = 14617
$ws4.Range("F10").Value = 170
$ws4.Range("F11").Value = 129
$ws4.Range("F12").Value = 5838
$ws4.Range("F13").Value = 595
$ws4.Range("F14").Value = 79
$ws4.Range("F15").Value = 45
$ws4.Range("F16").Value = 69
$ws4.Range("F18").Value = 16
$ws4.Range("F19").Value = 86
$ws4.Range("F20").Value = 186
$ws4.Range("F21").Value = 800
$ws4.Range("F23").Value = 61
$ws4.Range("F25").Value = 10623
$ws4.Range("F27").Value = 63
$ws4.Range("F28").Value = 90
$ws4.Range("F29").Value = 3739
$ws4.Range("F30").Value = 246
$ws4.Range("F31").Value = 68
